# Mission 6 voices in progress, weapons refactoring
#
# Adds a new voice-line row (DX_M06_0225_alaric) to the "m06" sheet,
# inserted before the existing "DX_M06_0230_alaric" row, and introduces a
# new column C that duplicates the "name" column (A) for every data row.
# Also turns on AutoFilter for the sheet's data range and updates the
# sheet's saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("m06")

# --- Insert the new row for the new voice line, right above the old row 25
#     (DX_M06_0230_alaric), shifting everything from row 25 down by one.
$ws.Rows.Item(25).Insert()

$ws.Cells.Item(25, 1).Value = "DX_M06_0225_alaric"
$ws.Cells.Item(25, 2).Value = "0xA8988D88"

# --- New column C: for every data row (2..47) duplicate the "name" value
#     that's already in column A.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 1).Value()
}

# --- Give the new column a sensible width (matches the widened column in
#     the authored workbook).
$ws.Columns.Item(3).ColumnWidth = 31.08

# --- Turn on AutoFilter across the full (now 47-row) table and register
#     the corresponding sheet-scoped _FilterDatabase name.
$ws.Range("A1:D47").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "='m06'!`$A`$1:`$D`$47")
try { $fdb.Visible = $false } catch { }

# --- Update the saved selection/view (author ended up with E17 selected,
#     scrolled back to the top of the sheet).
$ws.Range("E17").Select() | Out-Null
